# "on croyait que c'etait fini mais en fait non"
# Extend the "Occupation voie chantier" sheet with 40 additional rows
# (632-671), continuing the 15-minute timestamp series in column A with
# REC/FOR/DEP (columns B/C/D) all at 0, same as the tail of existing data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Occupation voie chantier")

# Re-use the existing date/time number format already applied to column A
# (same style as row 631) so the new timestamps render identically.
$fmt = $ws.Range("A631").NumberFormat

$lastRow = 631

# Horodatage values for the 40 new rows (15-minute steps continuing from
# row 631's 44787.55208333334), written as literal constants so the exact
# same double values as the source data are stored.
$timestamps = @(
    44787.5625,
    44787.57291666666,
    44787.58333333334,
    44787.59375,
    44787.60416666666,
    44787.61458333334,
    44787.625,
    44787.63541666666,
    44787.64583333334,
    44787.65625,
    44787.66666666666,
    44787.67708333334,
    44787.6875,
    44787.69791666666,
    44787.70833333334,
    44787.71875,
    44787.72916666666,
    44787.73958333334,
    44787.75,
    44787.76041666666,
    44787.77083333334,
    44787.78125,
    44787.79166666666,
    44787.80208333334,
    44787.8125,
    44787.82291666666,
    44787.83333333334,
    44787.84375,
    44787.85416666666,
    44787.86458333334,
    44787.875,
    44787.88541666666,
    44787.89583333334,
    44787.90625,
    44787.91666666666,
    44787.92708333334,
    44787.9375,
    44787.94791666666,
    44787.95833333334,
    44787.96875
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $r = $lastRow + 1 + $i

    $ws.Range("A$r").NumberFormat = $fmt
    $ws.Range("A$r").Value = $timestamps[$i]

    $ws.Range("B$r").Value = 0
    $ws.Range("C$r").Value = 0
    $ws.Range("D$r").Value = 0
}
